$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the old "Late" column (N), shifting
# the existing N/O/P columns to O/P/Q.
$ws.Columns("N").Insert() | Out-Null

# The sheet was left active with R8 selected when the workbook was saved.
$ws.Activate() | Out-Null
$ws.Range("R8").Select() | Out-Null
